# Rebuild the sheet data per the NATMI re-run (Dr Hou advice: include ECs as a sending cluster too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so the sheet (and shared-string table) rebuilds cleanly from the new data.
$ws.Cells.ClearContents()

# Header row (row 1)
$ws.Range("A1").Value = 'Sending cluster'
$ws.Range("B1").Value = 'Ligand symbol'
$ws.Range("C1").Value = 'Receptor symbol'
$ws.Range("D1").Value = 'Target cluster'
$ws.Range("E1").Value = 'Ligand-expressing cells'
$ws.Range("F1").Value = 'Ligand detection rate'
$ws.Range("G1").Value = 'Ligand average expression value'
$ws.Range("H1").Value = 'Ligand total expression value'
$ws.Range("I1").Value = 'Ligand derived specificity of average expression value'
$ws.Range("J1").Value = 'Ligand derived specificity of total expression value'
$ws.Range("K1").Value = 'Receptor-expressing cells'
$ws.Range("L1").Value = 'Receptor detection rate'
$ws.Range("M1").Value = 'Receptor average expression value'
$ws.Range("N1").Value = 'Receptor total expression value'
$ws.Range("O1").Value = 'Receptor derived specificity of average expression value'
$ws.Range("P1").Value = 'Receptor derived specificity of total expression value'
$ws.Range("Q1").Value = 'Edge average expression weight'
$ws.Range("R1").Value = 'Edge total expression weight'
$ws.Range("S1").Value = 'Edge average expression derived specificity'
$ws.Range("T1").Value = 'Edge total expression derived specificity'

# Data rows (rows 2-10): Sending cluster x {ECs, FAPs, sCs} against Inhba -> Acvr1b
# row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Inhba'
$ws.Range("C2").Value = 'Acvr1b'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.123204333333334
$ws.Range("H2").Value = 12.369613
$ws.Range("I2").Value = 0.2909967288544799
$ws.Range("J2").Value = 0.2909967288544799
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.453045666666667
$ws.Range("N2").Value = 13.359137
$ws.Range("O2").Value = 0.4394129038053478
$ws.Range("P2").Value = 0.4394129038053478
$ws.Range("Q2").Value = 18.36081718933123
$ws.Range("R2").Value = 165.247354703981
$ws.Range("S2").Value = 0.1278677176238045
$ws.Range("T2").Value = 0.1278677176238044

# row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Inhba'
$ws.Range("C3").Value = 'Acvr1b'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.123204333333334
$ws.Range("H3").Value = 12.369613
$ws.Range("I3").Value = 0.2909967288544799
$ws.Range("J3").Value = 0.2909967288544799
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.150099666666667
$ws.Range("N3").Value = 9.450299000000001
$ws.Range("O3").Value = 0.3108421842981904
$ws.Range("P3").Value = 0.3108421842981904
$ws.Range("Q3").Value = 12.98850459603189
$ws.Range("R3").Value = 116.896541364287
$ws.Range("S3").Value = 0.0904540588207548
$ws.Range("T3").Value = 0.09045405882075477

# row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Inhba'
$ws.Range("C4").Value = 'Acvr1b'
$ws.Range("D4").Value = 'sCs'
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.123204333333334
$ws.Range("H4").Value = 12.369613
$ws.Range("I4").Value = 0.2909967288544799
$ws.Range("J4").Value = 0.2909967288544799
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.530935
$ws.Range("N4").Value = 7.592805
$ws.Range("O4").Value = 0.2497449118964618
$ws.Range("P4").Value = 0.2497449118964618
$ws.Range("Q4").Value = 10.435562159385
$ws.Range("R4").Value = 93.92005943446502
$ws.Range("S4").Value = 0.07267495240992068
$ws.Range("T4").Value = 0.07267495240992068

# row 5
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Inhba'
$ws.Range("C5").Value = 'Acvr1b'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.433639666666666
$ws.Range("H5").Value = 25.300919
$ws.Range("I5").Value = 0.5952073574179045
$ws.Range("J5").Value = 0.5952073574179045
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.453045666666667
$ws.Range("N5").Value = 13.359137
$ws.Range("O5").Value = 0.4394129038053478
$ws.Range("P5").Value = 0.4394129038053478
$ws.Range("Q5").Value = 37.55538257187811
$ws.Range("R5").Value = 337.998443146903
$ws.Range("S5").Value = 0.2615417932893089
$ws.Range("T5").Value = 0.2615417932893089

# row 6
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Inhba'
$ws.Range("C6").Value = 'Acvr1b'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.433639666666666
$ws.Range("H6").Value = 25.300919
$ws.Range("I6").Value = 0.5952073574179045
$ws.Range("J6").Value = 0.5952073574179045
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.150099666666667
$ws.Range("N6").Value = 9.450299000000001
$ws.Range("O6").Value = 0.3108421842981904
$ws.Range("P6").Value = 0.3108421842981904
$ws.Range("Q6").Value = 26.56680550275345
$ws.Range("R6").Value = 239.101249524781
$ws.Range("S6").Value = 0.1850155550901352
$ws.Range("T6").Value = 0.1850155550901351

# row 7
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Inhba'
$ws.Range("C7").Value = 'Acvr1b'
$ws.Range("D7").Value = 'sCs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.433639666666666
$ws.Range("H7").Value = 25.300919
$ws.Range("I7").Value = 0.5952073574179045
$ws.Range("J7").Value = 0.5952073574179045
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.530935
$ws.Range("N7").Value = 7.592805
$ws.Range("O7").Value = 0.2497449118964618
$ws.Range("P7").Value = 0.2497449118964618
$ws.Range("Q7").Value = 21.344993809755
$ws.Range("R7").Value = 192.104944287795
$ws.Range("S7").Value = 0.1486500090384604
$ws.Range("T7").Value = 0.1486500090384604

# row 8
$ws.Range("A8").Value = 'sCs'
$ws.Range("B8").Value = 'Inhba'
$ws.Range("C8").Value = 'Acvr1b'
$ws.Range("D8").Value = 'ECs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.612402333333333
$ws.Range("H8").Value = 4.837207
$ws.Range("I8").Value = 0.1137959137276156
$ws.Range("J8").Value = 0.1137959137276156
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.453045666666667
$ws.Range("N8").Value = 13.359137
$ws.Range("O8").Value = 0.4394129038053478
$ws.Range("P8").Value = 0.4394129038053478
$ws.Range("Q8").Value = 7.180101223373224
$ws.Range("R8").Value = 64.62091101035901
$ws.Range("S8").Value = 0.05000339289223441
$ws.Range("T8").Value = 0.0500033928922344

# row 9
$ws.Range("A9").Value = 'sCs'
$ws.Range("B9").Value = 'Inhba'
$ws.Range("C9").Value = 'Acvr1b'
$ws.Range("D9").Value = 'FAPs'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.612402333333333
$ws.Range("H9").Value = 4.837207
$ws.Range("I9").Value = 0.1137959137276156
$ws.Range("J9").Value = 0.1137959137276156
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.150099666666667
$ws.Range("N9").Value = 9.450299000000001
$ws.Range("O9").Value = 0.3108421842981904
$ws.Range("P9").Value = 0.3108421842981904
$ws.Range("Q9").Value = 5.07922805276589
$ws.Range("R9").Value = 45.713052474893
$ws.Range("S9").Value = 0.03537257038730046
$ws.Range("T9").Value = 0.03537257038730045

# row 10
$ws.Range("A10").Value = 'sCs'
$ws.Range("B10").Value = 'Inhba'
$ws.Range("C10").Value = 'Acvr1b'
$ws.Range("D10").Value = 'sCs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.612402333333333
$ws.Range("H10").Value = 4.837207
$ws.Range("I10").Value = 0.1137959137276156
$ws.Range("J10").Value = 0.1137959137276156
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.530935
$ws.Range("N10").Value = 7.592805
$ws.Range("O10").Value = 0.2497449118964618
$ws.Range("P10").Value = 0.2497449118964618
$ws.Range("Q10").Value = 4.080885499515
$ws.Range("R10").Value = 36.727969495635
$ws.Range("S10").Value = 0.02841995044808073
$ws.Range("T10").Value = 0.02841995044808072

